$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C7").Value = 0.73387205387205301
$ws.Range("F7").Value = 0.74568926384783096
$ws.Range("I7").Value = 0.74739281244412603
$ws.Range("L7").Value = 0.74464934178146802

$ws.Range("G9").Select()
